$d = $word.ActiveDocument
$p = $d.Paragraphs(1)

# Add a 5-pt-spaced paragraph border (top/left/bottom/right) around the first paragraph
$borders = $p.Range.Borders
$borders.DistanceFromTop = 5
$borders.DistanceFromLeft = 5
$borders.DistanceFromBottom = 5
$borders.DistanceFromRight = 5

# Increase the left indent from 120 twips (6pt) to 225 twips (11.25pt)
$p.Format.LeftIndent = 11.25

# Replace the placeholder id text (which also removes the trailing space run)
$d.Content.Find.Execute("**ID__AFFARS_5334_topic_2__ID** ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "**ID__AFFARS_SUBPART_5334_2__ID**", 2)
